{"js": "// Remove the trailing \"Requisitos\" section (its Heading2 title paragraph\n// plus the following bulleted requirement paragraph) from the end of the\n// document body, right before the \"Bibliografia\" text and the section\n// break. The \"Bibliografia\" paragraph itself is left untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"Requisitos\" heading paragraph (exact, trimmed match so we\n// don't depend on trailing manual line breaks / whitespace quirks).\nlet requisitosIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"Requisitos\") {\n    requisitosIndex = i;\n  }\n}\n\nif (requisitosIndex !== -1) {\n  // Delete every paragraph from \"Requisitos\" through the end of the body\n  // (in this document that is just the heading itself and the single\n  // \"LOQ4240 ... (Requisito fraco)\" bullet paragraph that follows it).\n  for (let i = items.length - 1; i >= requisitosIndex; i--) {\n    items[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Requisitos\" section (its Heading2 title paragraph\n# plus the following bulleted requirement paragraph) from the end of the\n# document body. The preceding \"Bibliografia\" paragraph is left untouched.\n\n$d = $word.ActiveDocument\n\n$finder = $d.Content\n$found = $finder.Find.Execute(\"Requisitos\")\n\nif ($found) {\n    $headingPara = $finder.Paragraphs(1)\n    $rangeStart = $headingPara.Range.Start\n    $docEnd = $d.Content.End\n\n    $deleteRange = $d.Range($rangeStart, $docEnd)\n    $deleteRange.Delete()\n}\n"}
